$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-valued cells remain text (avoid Excel auto-numeric conversion)

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.185.56"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.28%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.843.17"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.58%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.78"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.18%  "

# Row 6
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07422"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.67%  "

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.08%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.82"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.19%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07714"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.68%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.825.66"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.48%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.005"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.23%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6754"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.09%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "86.04"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.01%  "

# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.39%  "

# Row 17
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.160.52"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.19%  "

# Row 18
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008304"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.48%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "228.42"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.54%  "

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.12%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.188"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.30%  "

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.01%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "160.44"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.31%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.687"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.91%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1403"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.03%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.99"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.77%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.507"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.49%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.171"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.52%  "

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.04%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.191"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.41%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05320"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.86%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7599"

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.878"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.55%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.31%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.681"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.28%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.330.92"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.54%  "

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.13%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.732"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.34%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9240"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.89%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.966"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.88%  "

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.17%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "103.48"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.01%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.07970"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +12.73%  "

# Row 45
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.978.80"

# Row 46
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5163"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.58%  "

# Row 47
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.771"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.37%  "

# Row 48
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000121"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.99%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "63.88"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.22%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.144"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.65%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05944"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.26%  "
